$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-16: full Clec11a -> Itga11 LR-pair table, recomputed after adding the
# "Resolving-Mac" target cluster (sender x target = 3 x 5 = 15 data rows).
$rows = @(
    ,('ECs', 'Clec11a', 'Itga11', 'ECs', "1", "0.3333333333333333", "0.4723746666666667", "1.417124", "0.02676815992729067", "0.02676815992729067", "3", "1", "0.5620333333333333", "1.6861", "0.01672125362771917", "0.01672125362771917", "0.2654903084888889", "2.3894127764", "0.0004475971912915759", "0.0004475971912915759")
    ,('ECs', 'Clec11a', 'Itga11', 'FAPs', "1", "0.3333333333333333", "0.4723746666666667", "1.417124", "0.02676815992729067", "0.02676815992729067", "3", "1", "32.710306", "98.13091799999999", "0.9731759495871608", "0.9731759495871608", "15.45151989331466", "139.063679039832", "0.02605012945594208", "0.02605012945594208")
    ,('ECs', 'Clec11a', 'Itga11', 'Inflammatory-Mac', "1", "0.3333333333333333", "0.4723746666666667", "1.417124", "0.02676815992729067", "0.02676815992729067", "2", "0.6666666666666666", "0.1711366666666666", "0.5134099999999999", "0.005091547847107109", "0.005091547847107109", "0.0808406258711111", "0.7275656328399999", "0.0001362913670488156", "0.0001362913670488156")
    ,('ECs', 'Clec11a', 'Itga11', 'MuSCs', "1", "0.3333333333333333", "0.4723746666666667", "1.417124", "0.02676815992729067", "0.02676815992729067", "3", "1", "0.133027", "0.399081", "0.003957733597653634", "0.003957733597653634", "0.06283858478266666", "0.5655472630439999", "0.0001059412458916039", "0.0001059412458916039")
    ,('ECs', 'Clec11a', 'Itga11', 'Resolving-Mac', "1", "0.3333333333333333", "0.4723746666666667", "1.417124", "0.02676815992729067", "0.02676815992729067", "1", "0.3333333333333333", "0.03541066666666667", "0.106232", "0.001053515340359328", "0.001053515340359328", "0.01672710186311111", "0.150543916768", "2.820066711659255E-05", "2.820066711659254E-05")
    ,('FAPs', 'Clec11a', 'Itga11', 'ECs', "3", "1", "17.041313", "51.123939", "0.9656838605972748", "0.9656838605972748", "3", "1", "0.5620333333333333", "1.6861", "0.01672125362771917", "0.01672125362771917", "9.577785949766666", "86.20007354789999", "0.01614744475724203", "0.01614744475724203")
    ,('FAPs', 'Clec11a', 'Itga11', 'FAPs', "3", "1", "17.041313", "51.123939", "0.9656838605972748", "0.9656838605972748", "3", "1", "32.710306", "98.13091799999999", "0.9731759495871608", "0.9731759495871608", "557.4265628717778", "5016.839065846002", "0.9397803080377484", "0.9397803080377484")
    ,('FAPs', 'Clec11a', 'Itga11', 'Inflammatory-Mac', "3", "1", "17.041313", "51.123939", "0.9656838605972748", "0.9656838605972748", "2", "0.6666666666666666", "0.1711366666666666", "0.5134099999999999", "0.005091547847107109", "0.005091547847107109", "2.916393502443333", "26.24754152198999", "0.004916825581410136", "0.004916825581410136")
    ,('FAPs', 'Clec11a', 'Itga11', 'MuSCs', "3", "1", "17.041313", "51.123939", "0.9656838605972748", "0.9656838605972748", "3", "1", "0.133027", "0.399081", "0.003957733597653634", "0.003957733597653634", "2.266954744451", "20.402592700059", "0.003821919459797703", "0.003821919459797703")
    ,('FAPs', 'Clec11a', 'Itga11', 'Resolving-Mac', "3", "1", "17.041313", "51.123939", "0.9656838605972748", "0.9656838605972748", "1", "0.3333333333333333", "0.03541066666666667", "0.106232", "0.001053515340359328", "0.001053515340359328", "0.6034442542053333", "5.430998287847999", "0.001017362761076648", "0.001017362761076647")
    ,('MuSCs', 'Clec11a', 'Itga11', 'ECs', "2", "0.6666666666666666", "0.1331983333333333", "0.399595", "0.007547979475434553", "0.007547979475434553", "3", "1", "0.5620333333333333", "1.6861", "0.01672125362771917", "0.01672125362771917", "0.07486190327777777", "0.6737571295", "0.0001262116791855598", "0.0001262116791855598")
    ,('MuSCs', 'Clec11a', 'Itga11', 'FAPs', "2", "0.6666666666666666", "0.1331983333333333", "0.399595", "0.007547979475434553", "0.007547979475434553", "3", "1", "32.710306", "98.13091799999999", "0.9731759495871608", "0.9731759495871608", "4.356958242023333", "39.21262417821", "0.007345512093470421", "0.007345512093470421")
    ,('MuSCs', 'Clec11a', 'Itga11', 'Inflammatory-Mac', "2", "0.6666666666666666", "0.1331983333333333", "0.399595", "0.007547979475434553", "0.007547979475434553", "2", "0.6666666666666666", "0.1711366666666666", "0.5134099999999999", "0.005091547847107109", "0.005091547847107109", "0.02279511877222222", "0.20515606895", "3.843089864815744E-05", "3.843089864815744E-05")
    ,('MuSCs', 'Clec11a', 'Itga11', 'MuSCs', "2", "0.6666666666666666", "0.1331983333333333", "0.399595", "0.007547979475434553", "0.007547979475434553", "3", "1", "0.133027", "0.399081", "0.003957733597653634", "0.003957733597653634", "0.01771897468833333", "0.159470772195", "2.987289196432738E-05", "2.987289196432738E-05")
    ,('MuSCs', 'Clec11a', 'Itga11', 'Resolving-Mac', "2", "0.6666666666666666", "0.1331983333333333", "0.399595", "0.007547979475434553", "0.007547979475434553", "1", "0.3333333333333333", "0.03541066666666667", "0.106232", "0.001053515340359328", "0.001053515340359328", "0.004716641782222222", "0.04244977604", "7.951912166087654E-06", "7.951912166087652E-06")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $rowVals = $rows[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $col = $c + 1
        $cell = $rowVals[$c]
        if ($col -le 4) {
            $ws.Cells.Item($r, $col).Value = $cell
        } else {
            $ws.Cells.Item($r, $col).Value = [double]$cell
        }
    }
}
